# Weekly update of "Fruta / hortaliza" prices.
# Existing rows get shifted to make room for the newest observation while
# keeping the same physical rows (2, 7, 8, 9, 10); rows 3-6 are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45141
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 8500
$ws.Range("L2").Value = 9000
$ws.Range("M2").Value = 8800
$ws.Range("P2").Value = 587

# Row 7
$ws.Range("D7").Value = 45119
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 20000
$ws.Range("L7").Value = 20000
$ws.Range("M7").Value = 20000
$ws.Range("P7").Value = 1333

# Row 8
$ws.Range("D8").Value = 44749
$ws.Range("J8").Value = 90
$ws.Range("K8").Value = 17000
$ws.Range("L8").Value = 18000
$ws.Range("M8").Value = 17556
$ws.Range("P8").Value = 1170

# Row 9
$ws.Range("D9").Value = 44750
$ws.Range("J9").Value = 140
$ws.Range("K9").Value = 19000
$ws.Range("L9").Value = 20000
$ws.Range("M9").Value = 19571
$ws.Range("P9").Value = 1305

# Row 10
$ws.Range("D10").Value = 45091
$ws.Range("J10").Value = 40
$ws.Range("K10").Value = 20000
$ws.Range("L10").Value = 22000
$ws.Range("M10").Value = 21000
$ws.Range("P10").Value = 1400
